$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.483.62"
$ws.Range("E2").Value = "  -0.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.640.38"
$ws.Range("E3").Value = "  -1.46%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.70"
$ws.Range("E5").Value = "  -1.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.28"
$ws.Range("E6").Value = "  -0.79%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  +3.90%  "

$ws.Range("E9").Value = "  +3.49%  "

$ws.Range("E10").Value = "  -0.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.75"
$ws.Range("E11").Value = "  -3.36%  "

$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.59"
$ws.Range("E13").Value = "  -2.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000194"
$ws.Range("E14").Value = "  -0.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.117.10"
$ws.Range("E15").Value = "  -1.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.368.61"
$ws.Range("E16").Value = "  -0.42%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.654.50"
$ws.Range("E17").Value = "  -1.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.54"
$ws.Range("E18").Value = "  -0.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.72"
$ws.Range("E19").Value = "  -2.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.41"
$ws.Range("E20").Value = "  -1.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "347.09"
$ws.Range("E21").Value = "  -1.25%  "

$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.84"
$ws.Range("E23").Value = "  -2.40%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000111"
$ws.Range("E24").Value = "  +0.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.51"
$ws.Range("E25").Value = "  -3.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.63"
$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.57"
$ws.Range("E27").Value = "  -3.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.163"
$ws.Range("E28").Value = "  -3.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "535.44"
$ws.Range("E30").Value = "  -0.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.83"
$ws.Range("E31").Value = "  -3.90%  "

$ws.Range("E32").Value = "  -3.37%  "

$ws.Range("E33").Value = "  -1.56%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.35"
$ws.Range("E34").Value = "  -3.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.37"
$ws.Range("E35").Value = "  -0.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.418"
$ws.Range("E36").Value = "  -1.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.19"
$ws.Range("E37").Value = "  -1.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "154.50"
$ws.Range("E39").Value = "  -3.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.90"
$ws.Range("E40").Value = "  -3.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "159.56"
$ws.Range("E42").Value = "  -3.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.04"
$ws.Range("E43").Value = "  -1.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.29"
$ws.Range("E44").Value = "  +2.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0600"
$ws.Range("E45").Value = "  -2.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.30"
$ws.Range("E46").Value = "  -3.57%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.632"
$ws.Range("E47").Value = "  -2.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0253"
$ws.Range("E48").Value = "  -3.86%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0995"
$ws.Range("E49").Value = "  -0.64%  "

$ws.Range("E50").Value = "  +7.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.51"
$ws.Range("E51").Value = "  -3.73%  "
